$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment Schedule")

# Insert a new column before O (shifting old O->P, P->Q), then move old column N's
# content/style into the newly inserted column, leaving N blank.
$ws.Range("O1:O14").Insert(-4161) | Out-Null   # xlShiftToRight
$ws.Range("N1:N14").Cut($ws.Range("O1:O14")) | Out-Null

# Activate the Repayment Schedule sheet and select R5 (matches final saved UI state).
$ws.Activate() | Out-Null
$ws.Range("R5").Select() | Out-Null

Write-Output "done"
